$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 65, shifting existing rows 65..103 down to 66..104
$ws.Rows.Item(65).Insert()

# Populate the new row 65 (copy of the previous row 65 content, with updated
# date / volume / price fields per the source data update)
$ws.Cells.Item(65, 1).Value = 4
$ws.Cells.Item(65, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(65, 3).Value = "Los Lagos"
$ws.Cells.Item(65, 4).Value = 44567
$ws.Cells.Item(65, 5).Value = 10
$ws.Cells.Item(65, 6).Value = 100112022
$ws.Cells.Item(65, 7).Value = "Arveja Verde"
$ws.Cells.Item(65, 8).Value = "Sin especificar"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 35
$ws.Cells.Item(65, 11).Value = 27000
$ws.Cells.Item(65, 12).Value = 27000
$ws.Cells.Item(65, 13).Value = 27000
$ws.Cells.Item(65, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(65, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(65, 16).Value = 1080
$ws.Cells.Item(65, 17).Value = 25
$ws.Cells.Item(65, 18).Value = "Hortaliza"
